{"js": "// Apply the edit described by the diff:\n// 1. Replace the text of the \"Remove mode switching from GhostController...\"\n//    paragraph with \"Make sure that ghost transition properly between Ghost\n//    states and ghost frightened modes.\"\n// 2. Remove the whole \"Make sure that ghosts are is switching between\n//    scatter and chase modes.\" paragraph.\n\nconst paras = context.document.body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nconst oldFirstText = \"Remove mode switching from GhostController and incorporate to ghost class.\";\nconst newFirstText = \"Make sure that ghost transition properly between Ghost states and ghost frightened modes.\";\nconst secondText = \"Make sure that ghosts are is switching between scatter and chase modes.\";\n\nlet firstPara = null;\nlet secondPara = null;\n\nfor (let i = 0; i < paras.items.length; i++) {\n  const text = paras.items[i].text;\n  if (text === oldFirstText) {\n    firstPara = paras.items[i];\n  } else if (text === secondText) {\n    secondPara = paras.items[i];\n  }\n}\n\nif (firstPara) {\n  // Replace the run text (keeps the paragraph, its bookmark, and its\n  // numbering/list formatting intact) by inserting the new text and\n  // clearing out the old content.\n  firstPara.insertText(newFirstText, Word.InsertLocation.replace);\n}\n\nif (secondPara) {\n  secondPara.delete();\n}\n\nawait context.sync();\n", "ps1": "# Apply the edit described by the diff:\n# 1. Replace the text of the \"Remove mode switching from GhostController...\"\n#    paragraph with \"Make sure that ghost transition properly between Ghost\n#    states and ghost frightened modes.\"\n# 2. Remove the whole \"Make sure that ghosts are is switching between\n#    scatter and chase modes.\" paragraph.\n\n$d = $word.ActiveDocument\n\n$oldFirstText = \"Remove mode switching from GhostController and incorporate to ghost class.\"\n$newFirstText = \"Make sure that ghost transition properly between Ghost states and ghost frightened modes.\"\n$secondText = \"Make sure that ghosts are is switching between scatter and chase modes.\"\n\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`n\", \"`x07\")\n\n    if ($text -eq $secondText) {\n        $p.Range.Delete()\n    }\n    elseif ($text -eq $oldFirstText) {\n        # Build a range covering just the paragraph's text (excluding the\n        # trailing paragraph mark) so the paragraph mark / bookmark that\n        # follows it is preserved.\n        $r = $p.Range\n        $r.MoveEnd(1, -1) | Out-Null\n        $r.Text = $newFirstText\n    }\n}\n"}
